$wb = $excel.ActiveWorkbook

# The workbook for this repo tracks 漫展 (convention) listing stats.
# "想去人数" (F column = interested count) and "最低票价" (G column = min price)
# were refreshed for the "展览" and "全部类型" sheets.

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 197

    $ws.Range("F3").Value = 764
    $ws.Range("G3").Value = 65

    $ws.Range("F6").Value = 4533

    $ws.Range("F8").Value = 365

    $ws.Range("F9").Value = 1301

    $ws.Range("F12").Value = 895

    $ws.Range("F14").Value = 499

    $ws.Range("F16").Value = 236
}
